$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.978.51"
$ws.Range("E2").Value = "  -5.11%  "

$ws.Range("D3").Value = "2.977.57"
$ws.Range("E3").Value = "  -6.71%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.20"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.10"
$ws.Range("E6").Value = "  -7.75%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "2.973.18"
$ws.Range("E8").Value = "  -6.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  -6.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.09"
$ws.Range("E11").Value = "  -3.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  -3.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -6.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.66"
$ws.Range("E14").Value = "  -6.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.119"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16").Value = "3.461.47"
$ws.Range("E16").Value = "  -6.74%  "

$ws.Range("D17").Value = "59.968.63"
$ws.Range("E17").Value = "  -5.00%  "

$ws.Range("D18").Value = "2.961.09"
$ws.Range("E18").Value = "  -6.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("E19").Value = "  -6.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.21"
$ws.Range("E20").Value = "  -6.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -7.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.661"
$ws.Range("E22").Value = "  -5.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.00"
$ws.Range("E23").Value = "  -8.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.75"
$ws.Range("E24").Value = "  -5.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.92"
$ws.Range("E25").Value = "  -4.81%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.52"
$ws.Range("E28").Value = "  -5.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  -7.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  -7.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -8.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.32"
$ws.Range("E32").Value = "  -7.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0931"
$ws.Range("E33").Value = "  -9.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("E34").Value = "  -9.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.950"
$ws.Range("E35").Value = "  -8.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.59"
$ws.Range("E36").Value = "  -4.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.59"
$ws.Range("E37").Value = "  -3.47%  "

$ws.Range("D38").Value = "0.0₃0659"
$ws.Range("E38").Value = "  -7.33%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0358"
$ws.Range("E39").Value = "  -7.94%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.90"
$ws.Range("E40").Value = "  -2.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.109"
$ws.Range("E41").Value = "  -3.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "381.54"
$ws.Range("E42").Value = "  -6.58%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").Value = "  -7.04%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.621.96"
$ws.Range("E44").Value = "  -6.67%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.236"
$ws.Range("E46").Value = "  -6.87%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.98"
$ws.Range("E47").Value = "  -7.05%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.38"
$ws.Range("E48").Value = "  -5.04%  "

$ws.Range("E49").Value = "  -4.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.38"
$ws.Range("E50").Value = "  -7.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.36"
$ws.Range("E51").Value = "  -11.00%  "
